# 1014-MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment-Loanproduct.xlsx
#
# The product name stored in cell B1 of both worksheets gets its space
# removed ("...TR-1-Late Repayment" -> "...TR-1-LateRepayment"), and the
# active/selected view moves from ProductLoanInput!A29 to
# ProductLoanOutput!B1.

$wb  = $excel.ActiveWorkbook
$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "1014-MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment"

# Correct the product name on both the input and the output sheet.
$wsInput.Range("B1").Value  = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Move the selection / active sheet: ProductLoanInput no longer keeps its
# old A29 selection (collapses back to B1), and ProductLoanOutput becomes
# the active tab with B1 selected.
$wsInput.Range("B1").Select()
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
